$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29 -> Art_028
$ws.Range("B29").Value = "https://archiveofourown.org/works/24329308"
$ws.Range("C29").Value = "Cap"
$ws.Range("D29").Value = "thecaptainspeaks"
$ws.Range("E29").Value = "AO3"

# Row 30 -> Art_029
$ws.Range("B30").Value = "https://archiveofourown.org/users/majesdane"
$ws.Range("C30").Value = "Maple"
$ws.Range("D30").Value = "majesdane"
$ws.Range("E30").Value = "AO3"

# Row 31 -> Art_030 (username entered before the link, matching the author's edit order)
$ws.Range("C31").Value = "Jamie"
$ws.Range("D31").Value = "isailonships"
$ws.Range("E31").Value = "AO3"
$ws.Range("B31").Value = "https://archiveofourown.org/users/ISailOnShips/pseuds/ISailOnShips"

# Update selection/view state to match the author's saved session
$ws.Range("B38").Select()
